$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 375.90475
$ws.Range("I33").Value = 310.3684
$ws.Range("K33").Value = 310.3684
$ws.Range("M33").Value = -81.36840000000001

$ws.Range("H92").Value = 723.6316
$ws.Range("I92").Value = 708.2778
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 708.2778
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 539.7222
$ws.Range("N92").Value = -3496

$ws.Range("H132").Value = 62785.055
$ws.Range("I132").Value = 7654.4116
$ws.Range("K132").Value = 22963.2348
$ws.Range("M132").Value = -20433.2348

$ws.Range("H137").Value = 11568.161
$ws.Range("I137").Value = 3437.75
$ws.Range("J137").Value = 14396.131
$ws.Range("K137").Value = 10313.25
$ws.Range("L137").Value = 43188.393
$ws.Range("M137").Value = -7763.25
$ws.Range("N137").Value = -48288.393

$ws.Range("H138").Value = 5958.5
$ws.Range("I138").Value = 6350.6665
$ws.Range("J138").Value = 5723.2
$ws.Range("K138").Value = 19051.9995
$ws.Range("L138").Value = 17169.6
$ws.Range("M138").Value = -13911.9995
$ws.Range("N138").Value = -27449.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1487207.8
$ws.Range("I32").Value = 675141.7
$ws.Range("K32").Value = 675141.7
$ws.Range("M32").Value = -674854.7

$ws.Range("H45").Value = 100085470
$ws.Range("I45").Value = 140780
$ws.Range("K45").Value = 140780
$ws.Range("M45").Value = -140403

$ws.Range("H61").Value = 3666.6875
$ws.Range("I61").Value = 2706.625
$ws.Range("J61").Value = 4626.75
$ws.Range("K61").Value = 2706.625
$ws.Range("L61").Value = 4626.75
$ws.Range("M61").Value = -2494.625
$ws.Range("N61").Value = -5050.75

$ws.Range("H102").Value = 83336410
$ws.Range("I102").Value = 100002740
$ws.Range("K102").Value = 100002740
$ws.Range("M102").Value = -100001118

$ws.Range("H132").Value = 3352.1667
$ws.Range("I132").Value = 2709.8096
$ws.Range("J132").Value = 4251.467
$ws.Range("K132").Value = 8129.4288
$ws.Range("L132").Value = 12754.401
$ws.Range("M132").Value = -5599.4288
$ws.Range("N132").Value = -17814.401

$ws.Range("H136").Value = 3666.6875
$ws.Range("I136").Value = 2706.625
$ws.Range("J136").Value = 4626.75
$ws.Range("K136").Value = 8119.875
$ws.Range("L136").Value = 13880.25
$ws.Range("M136").Value = -5569.875
$ws.Range("N136").Value = -18980.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 52090828
$ws.Range("I20").Value = 59532012
$ws.Range("K20").Value = 59532012
$ws.Range("M20").Value = -59531765

$ws.Range("H86").Value = 2917.9062
$ws.Range("I86").Value = 2587.2307
$ws.Range("J86").Value = 4350.8335
$ws.Range("K86").Value = 2587.2307
$ws.Range("L86").Value = 4350.8335
$ws.Range("M86").Value = -1464.2307
$ws.Range("N86").Value = -6596.8335

$ws.Range("H89").Value = 2917.9062
$ws.Range("I89").Value = 2587.2307
$ws.Range("J89").Value = 4350.8335
$ws.Range("K89").Value = 12936.1535
$ws.Range("L89").Value = 21754.1675
$ws.Range("M89").Value = -7320.1535
$ws.Range("N89").Value = -32986.1675

$ws.Range("H94").Value = 66667692
$ws.Range("I94").Value = 74075100
$ws.Range("J94").Value = 950
$ws.Range("K94").Value = 74075100
$ws.Range("L94").Value = 950
$ws.Range("M94").Value = -74074649
$ws.Range("N94").Value = -1852

$ws.Range("H137").Value = 77549
$ws.Range("J137").Value = 77549
$ws.Range("L137").Value = 77549
$ws.Range("N137").Value = -87749

$ws.Range("H140").Value = 59999
$ws.Range("J140").Value = 59999
$ws.Range("L140").Value = 59999
$ws.Range("N140").Value = -70359

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2454811.5
$ws.Range("I31").Value = 4504.1665
$ws.Range("J31").Value = 2781519
$ws.Range("K31").Value = 4504.1665
$ws.Range("L31").Value = 2781519
$ws.Range("M31").Value = -4209.1665
$ws.Range("N31").Value = -2782109

$ws.Range("H34").Value = 2454811.5
$ws.Range("I34").Value = 4504.1665
$ws.Range("J34").Value = 2781519
$ws.Range("K34").Value = 4504.1665
$ws.Range("L34").Value = 2781519
$ws.Range("M34").Value = -4302.1665
$ws.Range("N34").Value = -2781923

$ws.Range("H58").Value = 3571.2856
$ws.Range("I58").Value = 2666.3333
$ws.Range("K58").Value = 2666.3333
$ws.Range("M58").Value = -2463.3333

$ws.Range("H69").Value = 29812
$ws.Range("I69").Value = 24874
$ws.Range("J69").Value = 34750
$ws.Range("K69").Value = 24874
$ws.Range("L69").Value = 34750
$ws.Range("M69").Value = -24125
$ws.Range("N69").Value = -36248

$ws.Range("H72").Value = 29812
$ws.Range("I72").Value = 24874
$ws.Range("J72").Value = 34750
$ws.Range("K72").Value = 74622
$ws.Range("L72").Value = 104250
$ws.Range("M72").Value = -70878
$ws.Range("N72").Value = -111738

$ws.Range("H132").Value = 2984.2104
$ws.Range("I132").Value = 3065.889
$ws.Range("K132").Value = 9197.667000000001
$ws.Range("M132").Value = -6667.667000000001

$ws.Range("H134").Value = 3758.7334
$ws.Range("I134").Value = 4108.45
$ws.Range("J134").Value = 3059.3
$ws.Range("K134").Value = 12325.35
$ws.Range("L134").Value = 9177.900000000001
$ws.Range("M134").Value = -9790.349999999999
$ws.Range("N134").Value = -14247.9

$ws.Range("H136").Value = 3571.2856
$ws.Range("I136").Value = 2666.3333
$ws.Range("K136").Value = 7998.999899999999
$ws.Range("M136").Value = -5448.999899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1393330.9
$ws.Range("I68").Value = 3150.6956
$ws.Range("J68").Value = 2045864.5
$ws.Range("K68").Value = 9452.086800000001
$ws.Range("L68").Value = 6137593.5
$ws.Range("M68").Value = -8641.086800000001
$ws.Range("N68").Value = -6139215.5

$ws.Range("H71").Value = 1393330.9
$ws.Range("I71").Value = 3150.6956
$ws.Range("J71").Value = 2045864.5
$ws.Range("K71").Value = 28356.2604
$ws.Range("L71").Value = 18412780.5
$ws.Range("M71").Value = -24300.2604
$ws.Range("N71").Value = -18420892.5

$ws.Range("H98").Value = 510

$ws.Range("H134").Value = 3576.5
$ws.Range("I134").Value = 3576.5
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10729.5
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5659.5
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 10237
$ws.Range("I39").Value = 10237
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 10237
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -9705
$ws.Range("N39").ClearContents()

$ws.Range("H70").Value = 35719040
$ws.Range("I70").Value = 62502850
$ws.Range("J70").Value = 7291.6665
$ws.Range("K70").Value = 62502850
$ws.Range("L70").Value = 7291.6665
$ws.Range("M70").Value = -62502580
$ws.Range("N70").Value = -7831.6665

$ws.Range("H73").Value = 35719040
$ws.Range("I73").Value = 62502850
$ws.Range("J73").Value = 7291.6665
$ws.Range("K73").Value = 62502850
$ws.Range("L73").Value = 7291.6665
$ws.Range("M73").Value = -62501914
$ws.Range("N73").Value = -9163.666499999999

$ws.Range("H132").Value = 2029
$ws.Range("I132").Value = 2244.4
$ws.Range("J132").Value = 1670
$ws.Range("K132").Value = 6733.200000000001
$ws.Range("L132").Value = 5010
$ws.Range("M132").Value = -4203.200000000001
$ws.Range("N132").Value = -10070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2631.1428
$ws.Range("I7").Value = 2517
$ws.Range("J7").Value = 2745.2856
$ws.Range("K7").Value = 2517
$ws.Range("L7").Value = 2745.2856
$ws.Range("M7").Value = -2405
$ws.Range("N7").Value = -2969.2856

$ws.Range("H16").Value = 1836.6364
$ws.Range("I16").Value = 1167.1111
$ws.Range("J16").Value = 4849.5
$ws.Range("K16").Value = 1167.1111
$ws.Range("L16").Value = 4849.5
$ws.Range("M16").Value = -997.1111000000001
$ws.Range("N16").Value = -5189.5

$ws.Range("H40").Value = 58529.824
$ws.Range("I40").Value = 91806.21000000001
$ws.Range("J40").Value = 6766.5557
$ws.Range("K40").Value = 91806.21000000001
$ws.Range("L40").Value = 6766.5557
$ws.Range("M40").Value = -91670.21000000001
$ws.Range("N40").Value = -7038.5557

$ws.Range("H126").Value = 2631.1428
$ws.Range("I126").Value = 2517
$ws.Range("J126").Value = 2745.2856
$ws.Range("K126").Value = 7551
$ws.Range("L126").Value = 8235.856800000001
$ws.Range("M126").Value = -5081
$ws.Range("N126").Value = -13175.8568

$ws.Range("H132").Value = 10741.875
$ws.Range("I132").Value = 7129.4443
$ws.Range("K132").Value = 21388.3329
$ws.Range("M132").Value = -18858.3329

$ws.Range("H136").Value = 13001.556
$ws.Range("I136").Value = 7502.25
$ws.Range("J136").Value = 17401
$ws.Range("K136").Value = 22506.75
$ws.Range("L136").Value = 52203
$ws.Range("M136").Value = -19956.75
$ws.Range("N136").Value = -57303

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33849.715
$ws.Range("I54").Value = 27017
$ws.Range("J54").Value = 42960
$ws.Range("K54").Value = 27017
$ws.Range("L54").Value = 42960
$ws.Range("M54").Value = -26497
$ws.Range("N54").Value = -44000

$ws.Range("H81").Value = 7361.1113
$ws.Range("J81").Value = 5966.6665
$ws.Range("L81").Value = 11933.333
$ws.Range("N81").Value = -14055.333

$ws.Range("H84").Value = 7361.1113
$ws.Range("J84").Value = 5966.6665
$ws.Range("L84").Value = 59666.665
$ws.Range("N84").Value = -70274.66500000001

$ws.Range("H107").Value = 624.1539
$ws.Range("I107").Value = 544.6
$ws.Range("J107").Value = 889.3333
$ws.Range("K107").Value = 1633.8
$ws.Range("L107").Value = 2667.9999
$ws.Range("M107").Value = 286.1999999999998
$ws.Range("N107").Value = -6507.9999

$ws.Range("H124").Value = 480424.5
$ws.Range("J124").Value = 480424.5
$ws.Range("L124").Value = 480424.5
$ws.Range("N124").Value = -490244.5

$ws.Range("H126").Value = 6146.909
$ws.Range("I126").Value = 6560.6313
$ws.Range("J126").Value = 3526.6667
$ws.Range("K126").Value = 19681.8939
$ws.Range("L126").Value = 10580.0001
$ws.Range("M126").Value = -17211.8939
$ws.Range("N126").Value = -15520.0001

$ws.Range("H132").Value = 1507.5428
$ws.Range("I132").Value = 1521.4193
$ws.Range("K132").Value = 4564.257900000001
$ws.Range("M132").Value = -2034.257900000001
